# "render termino de version entregable"
#
# R5_usuario.xlsx tabulates, per responsible, date-validation stats.
# HoracioM's row (row 6) is refreshed with the figures from the final
# ("termino de version") render of the report:
#
#   n_total             (B6): 1.516 -> 1.513
#   n_error_fecha       (C6): 3     -> 0
#   pct_error_fecha     (D6): 0,20  -> 0,00
#   pct_regla_operativa (F6): 18,73 -> 18,77
#
# All columns in this sheet are stored as plain text (shared strings), so
# values that look numeric ("1.513", "0") are entered with a leading
# apostrophe to force Excel to keep them as text rather than re-parsing
# them as numbers; the comma-decimal values ("0,00", "18,77") are already
# unambiguous text in this locale and need no such hint.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "'1.513"
$ws.Range("C6").Value = "'0"
$ws.Range("D6").Value = "0,00"
$ws.Range("F6").Value = "18,77"
